$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 48; $r -le 97; $r++) {
    $ws.Cells.Item($r, 3).Value = "LIDO"
}
